$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run RU (Russia, column L) 1001; without crop
# Updated simulation results for rows 2-6 (columns B, C, L)

$ws.Range("B2").Value = 0.563185705845955
$ws.Range("L2").Value = 0.581687077490783

$ws.Range("B3").Value = 0.549025299914154
$ws.Range("L3").Value = 0.606993083582337

$ws.Range("B4").Value = 0.484505093491948
$ws.Range("C4").Value = 0.532496566334215
$ws.Range("L4").Value = 0.595431312960594

$ws.Range("B5").Value = 0.472605946730509
$ws.Range("L5").Value = 0.459773904184331

$ws.Range("B6").Value = 0.373830435892481
$ws.Range("L6").Value = 0.344282629685112
